$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns we touch keep their original plain-text
# representation instead of being auto-coerced into numbers by Excel.
$touchedCells = @( `
  "D2","E2","D3","E3","E4","D5","E5","D6","E6","D7","E7","E8","D9","E9", `
  "E10","D11","E11","E12","E13","D14","E14","D15","E15","D16","E16", `
  "D17","E17","E18","D19","E19","D20","E20","D21","E21","D22","E22", `
  "D23","E23","D24","E24","D25","E25","E26","E27","D28","E28","D29","E29", `
  "E30","D31","E31","D32","E32","E33","E34","D35","E35","D36","E36", `
  "D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","E42","E43", `
  "D44","E44","D45","E45","D46","E46","E47","D48","E48","D49","E49", `
  "D50","E50","D51","E51" `
)
foreach ($addr in $touchedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.232.36"
$ws.Range("E2").Value = "  -2.84%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.673.58"
$ws.Range("E3").Value = "  -3.41%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "680.02"
$ws.Range("E5").Value = "  -3.94%  "

# Row 6 - Solana
$ws.Range("D6").Value = "158.85"
$ws.Range("E6").Value = "  -7.24%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.670.43"
$ws.Range("E7").Value = "  -3.52%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.04%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -6.21%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -10.00%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "7.06"
$ws.Range("E11").Value = "  -5.83%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -9.83%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -7.88%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.293.20"
$ws.Range("E14").Value = "  -3.36%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "32.24"
$ws.Range("E15").Value = "  -11.12%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.683.94"
$ws.Range("E16").Value = "  -2.51%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.253.80"
$ws.Range("E17").Value = "  -2.85%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -1.39%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "15.75"
$ws.Range("E19").Value = "  -10.09%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  -11.13%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "470.66"
$ws.Range("E21").Value = "  -8.88%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "9.82"
$ws.Range("E22").Value = "  -5.92%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.645"
$ws.Range("E23").Value = "  -9.71%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "79.16"
$ws.Range("E24").Value = "  -5.86%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.821.63"
$ws.Range("E25").Value = "  -3.12%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -12.66%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "10.82"
$ws.Range("E28").Value = "  -14.58%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "9.03"
$ws.Range("E29").Value = "  -12.56%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -11.62%  "

# Row 31 - Fetch.AI
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  -15.28%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "6.62"
$ws.Range("E32").Value = "  -10.17%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -10.46%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("E34").Value = "  +0.25%  "

# Row 35 - EthereumClassic
$ws.Range("D35").Value = "26.50"
$ws.Range("E35").Value = "  -9.01%  "

# Row 36 - Kaspa
$ws.Range("D36").Value = "0.159"
$ws.Range("E36").Value = "  -8.51%  "

# Row 37 - Aptos
$ws.Range("D37").Value = "8.07"
$ws.Range("E37").Value = "  -12.85%  "

# Row 38 - Filecoin
$ws.Range("D38").Value = "6.02"
$ws.Range("E38").Value = "  -6.97%  "

# Row 39 - was Stacks, now USDe (rows 39/40 swapped order)
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40 - was USDe, now Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "2.23"
$ws.Range("E40").Value = "  -8.81%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0895"
$ws.Range("E41").Value = "  -11.49%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.07%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -7.09%  "

# Row 44 - Monero
$ws.Range("D44").Value = "165.35"
$ws.Range("E44").Value = "  -1.75%  "

# Row 45 - OKB
$ws.Range("D45").Value = "47.72"
$ws.Range("E45").Value = "  -4.75%  "

# Row 46 - dogwifhat
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  -17.37%  "

# Row 47 - ONDO
$ws.Range("E47").Value = "  -7.96%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "27.73"
$ws.Range("E48").Value = "  -3.91%  "

# Row 49 - was SuiNetwork, now FLOKI (rows 49/50 swapped order)
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000268"
$ws.Range("E49").Value = "  -12.63%  "

# Row 50 - was FLOKI, now SuiNetwork
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "1.07"
$ws.Range("E50").Value = "  -6.89%  "

# Row 51 - Cosmos
$ws.Range("D51").Value = "7.83"
$ws.Range("E51").Value = "  -9.07%  "
